# The presentation currently has its slide-master theme ("theme2.xml")
# carrying the "Integral" colour scheme, while the unused/secondary theme
# part ("theme1.xml", only wired to the notes master) carries the default
# "Office Theme" colour scheme. The authored change swaps these two theme
# colour schemes so the deck's applied (slide-visible) theme becomes the
# plain Office colours instead of Integral.
#
# The font scheme and format scheme are already identical between the two
# themes, so the only substantive difference is the 12 theme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). We recolour the
# presentation's live theme colour scheme (reached through any slide) to
# the target "Office Theme" palette, which rewrites the shared theme part
# used by every slide/layout/master.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette (RGB packed as VBA-style BGR integers: R + G*256 + B*65536)
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
